$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report generation date (row 5, column D)
$ws.Range("D5").Value = "2026.02.03 06:10"

# Append new trade history rows 195-214
$ws.Range("A195").Value = "2026.01.26 10:35:52"
$ws.Range("B195").Value = 705535092
$ws.Range("C195").Value = "AUDUSD"
$ws.Range("D195").Value = "buy"
$ws.Range("E195").Value = 0.9399999999999999
$ws.Range("F195").Value = 0.69163
$ws.Range("I195").Value = "2026.01.26 18:16:55"
$ws.Range("J195").Value = 0.6933
$ws.Range("K195").Value = 0
$ws.Range("L195").Value = 0
$ws.Range("M195").Value = 156.98
$ws.Range("N195").Value = "[sl 0.69330]"
$ws.Range("A196").Value = "2026.01.26 11:05:42"
$ws.Range("B196").Value = 705627973
$ws.Range("C196").Value = "AUDUSD"
$ws.Range("D196").Value = "buy"
$ws.Range("E196").Value = 0.97
$ws.Range("F196").Value = 0.69154
$ws.Range("I196").Value = "2026.01.26 18:16:55"
$ws.Range("J196").Value = 0.6933
$ws.Range("K196").Value = 0
$ws.Range("L196").Value = 0
$ws.Range("M196").Value = 170.72
$ws.Range("N196").Value = "[sl 0.69330]"
$ws.Range("A197").Value = "2026.01.26 10:03:58"
$ws.Range("B197").Value = 705432349
$ws.Range("C197").Value = "XAUUSD"
$ws.Range("D197").Value = "buy"
$ws.Range("E197").Value = 0.2
$ws.Range("F197").Value = 5081.72
$ws.Range("I197").Value = "2026.01.26 22:03:24"
$ws.Range("J197").Value = 5039.12
$ws.Range("K197").Value = 0
$ws.Range("L197").Value = 0
$ws.Range("M197").Value = -852
$ws.Range("N197").Value = "[sl 5039.12]"
$ws.Range("A198").Value = "2026.01.27 11:16:10"
$ws.Range("B198").Value = 711618843
$ws.Range("C198").Value = "AUDUSD"
$ws.Range("D198").Value = "sell"
$ws.Range("E198").Value = 1.21
$ws.Range("F198").Value = 0.69097
$ws.Range("I198").Value = "2026.01.27 12:52:16"
$ws.Range("J198").Value = 0.6923
$ws.Range("K198").Value = 0
$ws.Range("L198").Value = 0
$ws.Range("M198").Value = -160.93
$ws.Range("N198").Value = "[sl 0.69230]"
$ws.Range("A199").Value = "2026.01.27 11:00:42"
$ws.Range("B199").Value = 711565301
$ws.Range("C199").Value = "XAUUSD"
$ws.Range("D199").Value = "buy"
$ws.Range("E199").Value = 0.14
$ws.Range("F199").Value = 5087.07
$ws.Range("I199").Value = "2026.01.27 23:53:54"
$ws.Range("J199").Value = 5175.05
$ws.Range("K199").Value = 0
$ws.Range("L199").Value = 0
$ws.Range("M199").Value = 1231.72
$ws.Range("N199").Value = "[tp 5175.05]"
$ws.Range("A200").Value = "2026.01.28 11:24:38"
$ws.Range("B200").Value = 717477122
$ws.Range("C200").Value = "XAUUSD"
$ws.Range("D200").Value = "buy"
$ws.Range("E200").Value = 0.18
$ws.Range("F200").Value = 5292.2
$ws.Range("I200").Value = "2026.01.28 14:04:58"
$ws.Range("J200").Value = 5247.96
$ws.Range("K200").Value = 0
$ws.Range("L200").Value = 0
$ws.Range("M200").Value = -796.3200000000001
$ws.Range("N200").Value = "[sl 5247.96]"
$ws.Range("A201").Value = "2026.01.28 11:20:31"
$ws.Range("B201").Value = 717465475
$ws.Range("C201").Value = "AUDUSD"
$ws.Range("D201").Value = "buy"
$ws.Range("E201").Value = 0.74
$ws.Range("F201").Value = 0.70081
$ws.Range("I201").Value = "2026.01.28 18:08:24"
$ws.Range("J201").Value = 0.69809
$ws.Range("K201").Value = 0
$ws.Range("L201").Value = 0
$ws.Range("M201").Value = -201.28
$ws.Range("N201").Value = "[sl 0.69809]"
$ws.Range("A202").Value = "2026.01.28 11:08:50"
$ws.Range("B202").Value = 717422388
$ws.Range("C202").Value = "AUDUSD"
$ws.Range("D202").Value = "buy"
$ws.Range("E202").Value = 0.77
$ws.Range("F202").Value = 0.70029
$ws.Range("I202").Value = "2026.01.28 18:18:06"
$ws.Range("J202").Value = 0.6976599999999999
$ws.Range("K202").Value = 0
$ws.Range("L202").Value = 0
$ws.Range("M202").Value = -202.51
$ws.Range("N202").Value = "[sl 0.69766]"
$ws.Range("A203").Value = "2026.01.28 10:45:32"
$ws.Range("B203").Value = 717319946
$ws.Range("C203").Value = "BNBUSD"
$ws.Range("D203").Value = "buy"
$ws.Range("E203").Value = 52.93
$ws.Range("F203").Value = 904.66
$ws.Range("I203").Value = "2026.01.29 05:49:44"
$ws.Range("J203").Value = 893.0599999999999
$ws.Range("K203").Value = 0
$ws.Range("L203").Value = -411.53
$ws.Range("M203").Value = -6139.88
$ws.Range("N203").Value = "[sl 893.06]"
$ws.Range("A204").Value = "2026.01.29 11:08:58"
$ws.Range("B204").Value = 725798101
$ws.Range("C204").Value = "AUDUSD"
$ws.Range("D204").Value = "buy"
$ws.Range("E204").Value = 0.5600000000000001
$ws.Range("F204").Value = 0.70724
$ws.Range("I204").Value = "2026.01.29 13:36:24"
$ws.Range("J204").Value = 0.7042
$ws.Range("K204").Value = 0
$ws.Range("L204").Value = 0
$ws.Range("M204").Value = -170.24
$ws.Range("N204").Value = "[sl 0.70420]"
$ws.Range("A205").Value = "2026.01.29 11:02:33"
$ws.Range("B205").Value = 725735954
$ws.Range("C205").Value = "AUDUSD"
$ws.Range("D205").Value = "buy"
$ws.Range("E205").Value = 0.57
$ws.Range("F205").Value = 0.7071
$ws.Range("I205").Value = "2026.01.29 13:36:32"
$ws.Range("J205").Value = 0.70414
$ws.Range("K205").Value = 0
$ws.Range("L205").Value = 0
$ws.Range("M205").Value = -168.72
$ws.Range("N205").Value = "[sl 0.70414]"
$ws.Range("A206").Value = "2026.01.29 11:21:56"
$ws.Range("B206").Value = 725882678
$ws.Range("C206").Value = "XAUUSD"
$ws.Range("D206").Value = "buy"
$ws.Range("E206").Value = 0.07000000000000001
$ws.Range("F206").Value = 5538.16
$ws.Range("I206").Value = "2026.01.29 18:04:02"
$ws.Range("J206").Value = 5446.12
$ws.Range("K206").Value = 0
$ws.Range("L206").Value = 0
$ws.Range("M206").Value = -644.28
$ws.Range("N206").Value = "[sl 5446.12]"
$ws.Range("A207").Value = "2026.01.30 10:15:51"
$ws.Range("B207").Value = 735345595
$ws.Range("C207").Value = "BNBUSD"
$ws.Range("D207").Value = "sell"
$ws.Range("E207").Value = 17.53
$ws.Range("F207").Value = 845.34
$ws.Range("I207").Value = "2026.01.30 10:31:48"
$ws.Range("J207").Value = 844.84
$ws.Range("K207").Value = 0
$ws.Range("L207").Value = 0
$ws.Range("M207").Value = 87.65000000000001
$ws.Range("N207").Value = "[tp 844.84]"
$ws.Range("A208").Value = "2026.01.30 10:41:12"
$ws.Range("B208").Value = 735524843
$ws.Range("C208").Value = "BNBUSD"
$ws.Range("D208").Value = "sell"
$ws.Range("E208").Value = 17.56
$ws.Range("F208").Value = 843.39
$ws.Range("I208").Value = "2026.01.30 11:34:48"
$ws.Range("J208").Value = 842.89
$ws.Range("K208").Value = 0
$ws.Range("L208").Value = 0
$ws.Range("M208").Value = 87.8
$ws.Range("N208").Value = "[tp 842.89]"
$ws.Range("A209").Value = "2026.01.30 11:48:12"
$ws.Range("B209").Value = 736070754
$ws.Range("C209").Value = "XAUUSD"
$ws.Range("D209").Value = "sell"
$ws.Range("E209").Value = 0.05
$ws.Range("F209").Value = 5115.75
$ws.Range("I209").Value = "2026.01.30 12:57:22"
$ws.Range("J209").Value = 5012.95
$ws.Range("K209").Value = 0
$ws.Range("L209").Value = 0
$ws.Range("M209").Value = 514
$ws.Range("N209").Value = "[sl 5012.95]"
$ws.Range("A210").Value = "2026.01.30 18:23:34"
$ws.Range("B210").Value = 739230973
$ws.Range("C210").Value = "SP500"
$ws.Range("D210").Value = "buy"
$ws.Range("E210").Value = 42.2
$ws.Range("F210").Value = 6953.37
$ws.Range("I210").Value = "2026.01.30 19:21:17"
$ws.Range("J210").Value = 6929.82
$ws.Range("K210").Value = 0
$ws.Range("L210").Value = 0
$ws.Range("M210").Value = -993.8099999999999
$ws.Range("N210").Value = "[sl 6929.82]"
$ws.Range("A211").Value = "2026.02.01 10:00:57"
$ws.Range("B211").Value = 743412751
$ws.Range("C211").Value = "BNBUSD"
$ws.Range("D211").Value = "sell"
$ws.Range("E211").Value = 17.09
$ws.Range("F211").Value = 775.98
$ws.Range("I211").Value = "2026.02.01 10:19:10"
$ws.Range("J211").Value = 775.48
$ws.Range("K211").Value = 0
$ws.Range("L211").Value = 0
$ws.Range("M211").Value = 85.45
$ws.Range("N211").Value = "[tp 775.48]"
$ws.Range("A212").Value = "2026.02.01 10:39:02"
$ws.Range("B212").Value = 743437269
$ws.Range("C212").Value = "BNBUSD"
$ws.Range("D212").Value = "sell"
$ws.Range("E212").Value = 16.12
$ws.Range("F212").Value = 776.12
$ws.Range("I212").Value = "2026.02.01 13:37:41"
$ws.Range("J212").Value = 775.63
$ws.Range("K212").Value = 0
$ws.Range("L212").Value = 0
$ws.Range("M212").Value = 78.98
$ws.Range("N212").Value = "[tp 775.63]"
$ws.Range("A213").Value = "2026.02.02 11:02:10"
$ws.Range("B213").Value = 749086422
$ws.Range("C213").Value = "SP500"
$ws.Range("D213").Value = "sell"
$ws.Range("E213").Value = 46.6
$ws.Range("F213").Value = 6877.45
$ws.Range("I213").Value = "2026.02.02 11:02:12"
$ws.Range("J213").Value = 6878.58
$ws.Range("K213").Value = 0
$ws.Range("L213").Value = 0
$ws.Range("M213").Value = -52.66
$ws.Range("N213").Value = "[sl 6878.58]"
$ws.Range("A214").Value = "2026.02.02 11:04:08"
$ws.Range("B214").Value = 749106194
$ws.Range("C214").Value = "XAUUSD"
$ws.Range("D214").Value = "sell"
$ws.Range("E214").Value = 0.03
$ws.Range("F214").Value = 4578.79
$ws.Range("I214").Value = "2026.02.02 14:40:22"
$ws.Range("J214").Value = 4777.19
$ws.Range("K214").Value = 0
$ws.Range("L214").Value = 0
$ws.Range("M214").Value = -595.2
$ws.Range("N214").Value = "[sl 4777.19]"
